$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely - this shifts B:F left to A:E and removes the
# per-cell formatting (style "1" / border+bold) that used to live in column A.
$ws.Range("A:A").Delete()
